$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph: "AI Generated Revision:" pitch paragraph.
#
# Change 1: "...user-friendly websites and in 6 months, I earned my
#            ITF+ certification. I have hands-on experience..."
#        -> "...user-friendly websites. I completed my ITF+
#            certification, and I have hands-on experience..."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " and in 6 months, I earned my ITF+ certification. I have",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". I completed my ITF+ certification, and I have", 2)

# ------------------------------------------------------------------
# Change 2: "...Would you be available for a quick chat this week?"
#        -> "...Would you be available for a quick chat this week
#            to discuss more?"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "available for a quick chat this week?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "available for a quick chat this week to discuss more?", 2)
